$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Range Analysis: merge the per-row ratio formulas in column E
#    into a single shared formula group starting at E2.
# ------------------------------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Analysis")
$wsRange.Range("E2:E8").Formula = "=C2/B2"

# ------------------------------------------------------------------
# 2) Close off High Analysis: same treatment for column E.
# ------------------------------------------------------------------
$wsClose = $wb.Worksheets.Item("Close off High Analysis")
$wsClose.Range("E2:E8").Formula = "=C2/B2"

# ------------------------------------------------------------------
# 3) RSI Analysis: user clicked on column F header, selecting the
#    whole column.
# ------------------------------------------------------------------
$wsRsi = $wb.Worksheets.Item("RSI Analysis")
$wsRsi.Activate()
$wsRsi.Range("F1:F1048576").Select() | Out-Null

# ------------------------------------------------------------------
# 4) Add the new "H4RSI Analysis" worksheet after "Close off High
#    Analysis" (i.e. at the end of the tab strip) with the H4RSI
#    band breakdown, mirroring the layout of "RSI Analysis".
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNew.Name = "H4RSI Analysis"

$wsNew.Range("A1").Value = "RSIBand"
$wsNew.Range("B1").Value = "(No column name)"
$wsNew.Range("C1").Value = "Trades"
$wsNew.Range("D1").Value = "Winners"
$wsNew.Range("E1").Value = "Losers"
$wsNew.Range("F1").Value = "WL Ratio"
$wsNew.Range("F1").NumberFormat = "0%"

$wsNew.Range("A2").Value = "Below 55"
$wsNew.Range("B2").Value = 1
$wsNew.Range("C2").Value = 44238
$wsNew.Range("D2").Value = 12810
$wsNew.Range("E2").Value = 31428

$wsNew.Range("A3").Value = "Below 60"
$wsNew.Range("B3").Value = 2
$wsNew.Range("C3").Value = 125084
$wsNew.Range("D3").Value = 38341
$wsNew.Range("E3").Value = 86743

$wsNew.Range("A4").Value = "Below 65"
$wsNew.Range("B4").Value = 3
$wsNew.Range("C4").Value = 96043
$wsNew.Range("D4").Value = 27301
$wsNew.Range("E4").Value = 68742

$wsNew.Range("A5").Value = "Below 70"
$wsNew.Range("B5").Value = 4
$wsNew.Range("C5").Value = 43114
$wsNew.Range("D5").Value = 14068
$wsNew.Range("E5").Value = 29046

$wsNew.Range("A6").Value = "Below 75"
$wsNew.Range("B6").Value = 5
$wsNew.Range("C6").Value = 13871
$wsNew.Range("D6").Value = 4953
$wsNew.Range("E6").Value = 8918

$wsNew.Range("A7").Value = "Above 75"
$wsNew.Range("B7").Value = 6
$wsNew.Range("C7").Value = 5013
$wsNew.Range("D7").Value = 2074
$wsNew.Range("E7").Value = 2939

# F2 is entered on its own, then filled down F3:F7 as a shared group,
# exactly mirroring how the sibling "RSI Analysis" sheet is built.
$wsNew.Range("F2").Formula = "=D2/C2"
$wsNew.Range("F2").NumberFormat = "0%"
$wsNew.Range("F3:F7").Formula = "=D3/C3"
$wsNew.Range("F3:F7").NumberFormat = "0%"

$wsNew.Columns.Item(6).ColumnWidth = 8

$wsNew.Range("C14").Select() | Out-Null

# ------------------------------------------------------------------
# 5) Make "H4 Analysis 1" the active tab again (it was tab 0, the
#    default, so the explicit activeTab marker on the workbook goes
#    away and tabSelected moves off "Close off High Analysis").
# ------------------------------------------------------------------
$wsFirst = $wb.Worksheets.Item("H4 Analysis 1")
$wsFirst.Activate()

Write-Host "Added H4RSI Analysis sheet and refreshed formulas"
